$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-26 02:48:44'
$ws.Range('O2').Value = '1.7 °C'
$ws.Range('E3').Value = '2026-02-26 02:48:46'
$ws.Range('N3').Value = '0.5 °C 2:12 TU'
$ws.Range('O3').Value = '1.7 °C'
$ws.Range('E4').Value = '2026-02-26 02:48:48'
$ws.Range('N4').Value = '6.3 °C 2:23 TU'
$ws.Range('O4').Value = '8.7 °C'
$ws.Range('E5').Value = '2026-02-26 02:48:51'
$ws.Range('H5').NumberFormat = '@'
$ws.Range('H5').Value = '50%'
$ws.Range('N5').Value = '2.8 °C 2:15 TU'
$ws.Range('O5').Value = '3.5 °C'
$ws.Range('E6').Value = '2026-02-26 02:48:53'
$ws.Range('N6').Value = '8.7 °C 2:29 TU'
$ws.Range('O6').Value = '10.2 °C'
$ws.Range('E7').Value = '2026-02-26 02:48:56'
$ws.Range('N7').Value = '11.5 °C 2:21 TU'
$ws.Range('O7').Value = '11.9 °C'
$ws.Range('E8').Value = '2026-02-26 02:48:59'
$ws.Range('H8').NumberFormat = '@'
$ws.Range('H8').Value = '96%'
$ws.Range('E9').Value = '2026-02-26 02:49:01'
$ws.Range('L9').Value = '2.5 km/h - 124º 2:28 TU'
$ws.Range('N9').Value = '10.8 °C 2:24 TU'
$ws.Range('E10').Value = '2026-02-26 02:49:04'
$ws.Range('N10').Value = '3.9 °C 2:18 TU'
$ws.Range('O10').Value = '4.6 °C'
$ws.Range('E11').Value = '2026-02-26 02:49:06'
$ws.Range('H11').NumberFormat = '@'
$ws.Range('H11').Value = '92%'
$ws.Range('O11').Value = '2.4 °C'
$ws.Range('E12').Value = '2026-02-26 02:49:09'
$ws.Range('N12').Value = '8.6 °C 2:25 TU'
$ws.Range('O12').Value = '10.0 °C'
$ws.Range('E13').Value = '2026-02-26 02:49:11'
$ws.Range('H13').NumberFormat = '@'
$ws.Range('H13').Value = '92%'
$ws.Range('J13').Value = '1031.3 hPa'
$ws.Range('N13').Value = '-1.6 °C 2:25 TU'
$ws.Range('O13').Value = '-0.3 °C'
$ws.Range('E14').Value = '2026-02-26 02:49:14'
$ws.Range('L14').Value = '14.8 km/h - 327º 2:07 TU'
$ws.Range('O14').Value = '10.4 °C'
$ws.Range('E15').Value = '2026-02-26 02:49:17'
$ws.Range('N15').Value = '10.1 °C 2:28 TU'
$ws.Range('O15').Value = '10.9 °C'
$ws.Range('E16').Value = '2026-02-26 02:49:19'
$ws.Range('N16').Value = '1.2 °C 2:29 TU'
$ws.Range('E17').Value = '2026-02-26 02:49:22'
$ws.Range('E18').Value = '2026-02-26 02:49:24'
$ws.Range('J18').Value = '1026.3 hPa'
$ws.Range('L18').Value = '1.4 km/h - 298º 2:24 TU'
$ws.Range('N18').Value = '8.2 °C 2:25 TU'
$ws.Range('O18').Value = '8.7 °C'
$ws.Range('E19').Value = '2026-02-26 02:49:27'
$ws.Range('H19').NumberFormat = '@'
$ws.Range('H19').Value = '71%'
$ws.Range('E20').Value = '2026-02-26 02:49:30'
$ws.Range('H20').NumberFormat = '@'
$ws.Range('H20').Value = '52%'
$ws.Range('L20').Value = '18.0 km/h - 281º 2:12 TU'
$ws.Range('N20').Value = '-0.7 °C 2:06 TU'
$ws.Range('O20').Value = '1.2 °C'
$ws.Range('E21').Value = '2026-02-26 02:49:32'
$ws.Range('H21').NumberFormat = '@'
$ws.Range('H21').Value = '82%'
$ws.Range('J21').Value = '1028.3 hPa'
$ws.Range('N21').Value = '3.8 °C 2:10 TU'
$ws.Range('O21').Value = '4.9 °C'
$ws.Range('E22').Value = '2026-02-26 02:49:35'
$ws.Range('E23').Value = '2026-02-26 02:49:38'
$ws.Range('H23').NumberFormat = '@'
$ws.Range('H23').Value = '45%'
$ws.Range('N23').Value = '1.8 °C 2:23 TU'
$ws.Range('O23').Value = '2.7 °C'
$ws.Range('E24').Value = '2026-02-26 02:49:40'
$ws.Range('N24').Value = '8.1 °C 2:27 TU'
$ws.Range('O24').Value = '8.8 °C'
$ws.Range('E25').Value = '2026-02-26 02:49:42'
$ws.Range('H25').NumberFormat = '@'
$ws.Range('H25').Value = '41%'
$ws.Range('L25').Value = '23.0 km/h - 6º 2:16 TU'
$ws.Range('O25').Value = '3.0 °C'
$ws.Range('E26').Value = '2026-02-26 02:49:45'
$ws.Range('H26').NumberFormat = '@'
$ws.Range('H26').Value = '50%'
$ws.Range('J26').Value = '1025.1 hPa'
$ws.Range('M26').Value = '8.3 °C 2:20 TU'
$ws.Range('O26').Value = '6.7 °C'
$ws.Range('E27').Value = '2026-02-26 02:49:48'
$ws.Range('N27').Value = '1.9 °C 2:14 TU'
$ws.Range('O27').Value = '2.7 °C'
$ws.Range('E28').Value = '2026-02-26 02:49:51'
$ws.Range('N28').Value = '8.1 °C 2:29 TU'
$ws.Range('O28').Value = '9.0 °C'
$ws.Range('E29').Value = '2026-02-26 02:49:53'
$ws.Range('E30').Value = '2026-02-26 02:49:56'
$ws.Range('N30').Value = '10.8 °C 2:18 TU'
$ws.Range('E31').Value = '2026-02-26 02:49:58'
$ws.Range('N31').Value = '10.2 °C 2:29 TU'
$ws.Range('O31').Value = '10.8 °C'
$ws.Range('E32').Value = '2026-02-26 02:50:01'
$ws.Range('O32').Value = '1.8 °C'
$ws.Range('E33').Value = '2026-02-26 02:50:04'
$ws.Range('H33').NumberFormat = '@'
$ws.Range('H33').Value = '74%'
$ws.Range('N33').Value = '2.4 °C 2:29 TU'
$ws.Range('O33').Value = '3.2 °C'
$ws.Range('E34').Value = '2026-02-26 02:50:06'
$ws.Range('H34').NumberFormat = '@'
$ws.Range('H34').Value = '59%'
$ws.Range('L34').Value = '20.2 km/h - 36º 2:23 TU'
$ws.Range('M34').Value = '5.6 °C 2:21 TU'
$ws.Range('O34').Value = '0.6 °C'
$ws.Range('E35').Value = '2026-02-26 02:50:09'
$ws.Range('H35').NumberFormat = '@'
$ws.Range('H35').Value = '28%'
$ws.Range('E36').Value = '2026-02-26 02:50:12'
$ws.Range('E37').Value = '2026-02-26 02:50:14'
$ws.Range('H37').NumberFormat = '@'
$ws.Range('H37').Value = '90%'
$ws.Range('L37').Value = '16.2 km/h - 238º 2:27 TU'
$ws.Range('E38').Value = '2026-02-26 02:50:17'
$ws.Range('N38').Value = '7.2 °C 2:29 TU'
$ws.Range('O38').Value = '9.1 °C'
$ws.Range('E39').Value = '2026-02-26 02:50:19'
$ws.Range('H39').NumberFormat = '@'
$ws.Range('H39').Value = '39%'
$ws.Range('L39').Value = '32.0 km/h - 327º 2:05 TU'
$ws.Range('O39').Value = '2.9 °C'
$ws.Range('E40').Value = '2026-02-26 02:50:22'
$ws.Range('J40').Value = '1029.1 hPa'
$ws.Range('N40').Value = '2.1 °C 2:07 TU'
$ws.Range('O40').Value = '2.7 °C'
$ws.Range('E41').Value = '2026-02-26 02:50:24'
$ws.Range('J41').Value = '1025.6 hPa'
$ws.Range('N41').Value = '6.9 °C 2:17 TU'
$ws.Range('O41').Value = '8.6 °C'
$ws.Range('E42').Value = '2026-02-26 02:50:27'
$ws.Range('O42').Value = '8.4 °C'
$ws.Range('E43').Value = '2026-02-26 02:50:29'
$ws.Range('H43').NumberFormat = '@'
$ws.Range('H43').Value = '94%'
$ws.Range('O43').Value = '3.6 °C'
$ws.Range('E44').Value = '2026-02-26 02:50:32'
$ws.Range('N44').Value = '-1.0 °C 2:07 TU'
$ws.Range('E45').Value = '2026-02-26 02:50:34'
$ws.Range('H45').NumberFormat = '@'
$ws.Range('H45').Value = '58%'
$ws.Range('J45').Value = '1026.9 hPa'
$ws.Range('N45').Value = '4.6 °C 2:29 TU'
$ws.Range('O45').Value = '6.6 °C'
$ws.Range('E46').Value = '2026-02-26 02:50:37'
$ws.Range('H46').NumberFormat = '@'
$ws.Range('H46').Value = '95%'
$ws.Range('N46').Value = '7.1 °C 2:02 TU'
$ws.Range('O46').Value = '8.4 °C'
